$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.650236010551453
$ws.Range("B1").Value = 1.658631682395935
$ws.Range("C1").Value = 1.785469651222229
$ws.Range("D1").Value = 2.554950952529907
$ws.Range("E1").Value = 2.865405559539795
